$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Clear the "Nationalität" value for row 17 (Paul) - C17
$ws.Range("C17").ClearContents()

# Update selection to D25
$ws.Range("D25").Select()
